$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-13 are being permuted (entire row contents move to different row numbers).
# Build the full replacement block as a 2D array (rows 3..13, columns A..AY) and assign it in one shot.
$data = New-Object "object[,]" 11,51

# Row 3
$data[0,0] = 101185517
$data[0,1] = 89832
$data[0,2] = 'Ovaliderad'
$data[0,3] = 'VU'
$data[0,4] = 1209
$data[0,5] = 'Rynkskinn'
$data[0,6] = 'Phlebia centrifuga'
$data[0,7] = 'P.Karst.'
$data[0,8] = ""
$data[0,9] = ""
$data[0,10] = ""
$data[0,11] = ""
$data[0,12] = ""
$data[0,13] = ""
$data[0,14] = ""
$data[0,15] = 'Knätten NB, Hjd'
$data[0,16] = 443198.490593656
$data[0,17] = 6909846.427368461
$data[0,18] = 25
$data[0,19] = 'Jämtland'
$data[0,20] = 'Härjedalen'
$data[0,21] = 'Härjedalen'
$data[0,22] = 'Vemdalen'
$data[0,23] = ""
$data[0,24] = "'2022-05-27"
$data[0,25] = "'00:00"
$data[0,26] = "'2022-05-27"
$data[0,27] = "'00:00"
$data[0,28] = ""
$data[0,29] = $false
$data[0,30] = $false
$data[0,31] = ""
$data[0,32] = $false
$data[0,33] = ""
$data[0,34] = ""
$data[0,35] = ""
$data[0,36] = ""
$data[0,37] = ""
$data[0,38] = ""
$data[0,39] = ""
$data[0,40] = ""
$data[0,41] = ""
$data[0,42] = ""
$data[0,43] = ""
$data[0,44] = ""
$data[0,45] = ""
$data[0,46] = ""
$data[0,47] = ""
$data[0,48] = 'Alexander Singer'
$data[0,49] = 'Alexander Singer'
$data[0,50] = ""

# Row 4
$data[1,0] = 101185534
$data[1,1] = 89392
$data[1,2] = 'Ovaliderad'
$data[1,3] = 'NT'
$data[1,4] = 1202
$data[1,5] = 'Ullticka'
$data[1,6] = 'Phellinidium ferrugineofuscum'
$data[1,7] = '(P.Karst.) Fiasson & Niemelä'
$data[1,8] = ""
$data[1,9] = ""
$data[1,10] = ""
$data[1,11] = ""
$data[1,12] = ""
$data[1,13] = ""
$data[1,14] = ""
$data[1,15] = 'Knätten NB, Hjd'
$data[1,16] = 443198.490593656
$data[1,17] = 6909846.427368461
$data[1,18] = 25
$data[1,19] = 'Jämtland'
$data[1,20] = 'Härjedalen'
$data[1,21] = 'Härjedalen'
$data[1,22] = 'Vemdalen'
$data[1,23] = ""
$data[1,24] = "'2022-05-27"
$data[1,25] = "'00:00"
$data[1,26] = "'2022-05-27"
$data[1,27] = "'00:00"
$data[1,28] = ""
$data[1,29] = $false
$data[1,30] = $false
$data[1,31] = ""
$data[1,32] = $false
$data[1,33] = ""
$data[1,34] = ""
$data[1,35] = ""
$data[1,36] = ""
$data[1,37] = ""
$data[1,38] = ""
$data[1,39] = ""
$data[1,40] = ""
$data[1,41] = ""
$data[1,42] = ""
$data[1,43] = ""
$data[1,44] = ""
$data[1,45] = ""
$data[1,46] = ""
$data[1,47] = ""
$data[1,48] = 'Alexander Singer'
$data[1,49] = 'Alexander Singer'
$data[1,50] = ""

# Row 5
$data[2,0] = 111221736
$data[2,1] = 56398
$data[2,2] = 'Ovaliderad'
$data[2,3] = 'NT'
$data[2,4] = 100109
$data[2,5] = 'Tretåig hackspett'
$data[2,6] = 'Picoides tridactylus'
$data[2,7] = '(Linnaeus, 1758)'
$data[2,8] = ""
$data[2,9] = ""
$data[2,10] = ""
$data[2,11] = ""
$data[2,12] = ""
$data[2,13] = ""
$data[2,14] = ""
$data[2,15] = 'A 25562-2023, Hjd'
$data[2,16] = 443249.6264723797
$data[2,17] = 6909840.911127058
$data[2,18] = 10
$data[2,19] = 'Jämtland'
$data[2,20] = 'Härjedalen'
$data[2,21] = 'Härjedalen'
$data[2,22] = 'Vemdalen'
$data[2,23] = ""
$data[2,24] = "'2023-07-21"
$data[2,25] = "'00:00"
$data[2,26] = "'2023-07-21"
$data[2,27] = "'00:00"
$data[2,28] = 'Ringhack'
$data[2,29] = $false
$data[2,30] = $false
$data[2,31] = ""
$data[2,32] = $false
$data[2,33] = ""
$data[2,34] = ""
$data[2,35] = ""
$data[2,36] = ""
$data[2,37] = ""
$data[2,38] = ""
$data[2,39] = ""
$data[2,40] = ""
$data[2,41] = ""
$data[2,42] = ""
$data[2,43] = ""
$data[2,44] = ""
$data[2,45] = ""
$data[2,46] = ""
$data[2,47] = ""
$data[2,48] = 'Christer Johansson'
$data[2,49] = 'Christer Johansson'
$data[2,50] = ""

# Row 6
$data[3,0] = 111221709
$data[3,1] = 56398
$data[3,2] = 'Ovaliderad'
$data[3,3] = 'NT'
$data[3,4] = 100109
$data[3,5] = 'Tretåig hackspett'
$data[3,6] = 'Picoides tridactylus'
$data[3,7] = '(Linnaeus, 1758)'
$data[3,8] = ""
$data[3,9] = ""
$data[3,10] = ""
$data[3,11] = ""
$data[3,12] = ""
$data[3,13] = ""
$data[3,14] = ""
$data[3,15] = 'A 25562-2023, Hjd'
$data[3,16] = 443254.9775056695
$data[3,17] = 6909826.869210822
$data[3,18] = 10
$data[3,19] = 'Jämtland'
$data[3,20] = 'Härjedalen'
$data[3,21] = 'Härjedalen'
$data[3,22] = 'Vemdalen'
$data[3,23] = ""
$data[3,24] = "'2023-07-21"
$data[3,25] = "'00:00"
$data[3,26] = "'2023-07-21"
$data[3,27] = "'00:00"
$data[3,28] = 'Ringhack'
$data[3,29] = $false
$data[3,30] = $false
$data[3,31] = ""
$data[3,32] = $false
$data[3,33] = ""
$data[3,34] = ""
$data[3,35] = ""
$data[3,36] = ""
$data[3,37] = ""
$data[3,38] = ""
$data[3,39] = ""
$data[3,40] = ""
$data[3,41] = ""
$data[3,42] = ""
$data[3,43] = ""
$data[3,44] = ""
$data[3,45] = ""
$data[3,46] = ""
$data[3,47] = ""
$data[3,48] = 'Christer Johansson'
$data[3,49] = 'Christer Johansson'
$data[3,50] = ""

# Row 7
$data[4,0] = 111221699
$data[4,1] = 56543
$data[4,2] = 'Ovaliderad'
$data[4,3] = 'NT'
$data[4,4] = 103021
$data[4,5] = 'Talltita'
$data[4,6] = 'Poecile montanus'
$data[4,7] = '(Conrad von Baldenstein, 1827)'
$data[4,8] = "'1"
$data[4,9] = ""
$data[4,10] = ""
$data[4,11] = ""
$data[4,12] = 'spel/sång'
$data[4,13] = ""
$data[4,14] = ""
$data[4,15] = 'A 25562-2023, Hjd'
$data[4,16] = 443097.6233577073
$data[4,17] = 6909995.088246249
$data[4,18] = 10
$data[4,19] = 'Jämtland'
$data[4,20] = 'Härjedalen'
$data[4,21] = 'Härjedalen'
$data[4,22] = 'Vemdalen'
$data[4,23] = ""
$data[4,24] = "'2023-07-21"
$data[4,25] = "'00:00"
$data[4,26] = "'2023-07-21"
$data[4,27] = "'00:00"
$data[4,28] = ""
$data[4,29] = $false
$data[4,30] = $false
$data[4,31] = ""
$data[4,32] = $false
$data[4,33] = ""
$data[4,34] = ""
$data[4,35] = ""
$data[4,36] = ""
$data[4,37] = ""
$data[4,38] = ""
$data[4,39] = ""
$data[4,40] = ""
$data[4,41] = ""
$data[4,42] = ""
$data[4,43] = ""
$data[4,44] = ""
$data[4,45] = ""
$data[4,46] = ""
$data[4,47] = ""
$data[4,48] = 'Christer Johansson'
$data[4,49] = 'Christer Johansson'
$data[4,50] = ""

# Row 8
$data[5,0] = 101188502
$data[5,1] = 89392
$data[5,2] = 'Ovaliderad'
$data[5,3] = 'NT'
$data[5,4] = 1202
$data[5,5] = 'Ullticka'
$data[5,6] = 'Phellinidium ferrugineofuscum'
$data[5,7] = '(P.Karst.) Fiasson & Niemelä'
$data[5,8] = ""
$data[5,9] = ""
$data[5,10] = ""
$data[5,11] = ""
$data[5,12] = ""
$data[5,13] = ""
$data[5,14] = ""
$data[5,15] = 'Knätten NB, Hjd'
$data[5,16] = 443372.9596193716
$data[5,17] = 6909394.679763782
$data[5,18] = 25
$data[5,19] = 'Jämtland'
$data[5,20] = 'Härjedalen'
$data[5,21] = 'Härjedalen'
$data[5,22] = 'Vemdalen'
$data[5,23] = ""
$data[5,24] = "'2022-05-27"
$data[5,25] = "'00:00"
$data[5,26] = "'2022-05-27"
$data[5,27] = "'00:00"
$data[5,28] = ""
$data[5,29] = $false
$data[5,30] = $false
$data[5,31] = ""
$data[5,32] = $false
$data[5,33] = ""
$data[5,34] = ""
$data[5,35] = ""
$data[5,36] = ""
$data[5,37] = ""
$data[5,38] = ""
$data[5,39] = ""
$data[5,40] = ""
$data[5,41] = ""
$data[5,42] = ""
$data[5,43] = ""
$data[5,44] = ""
$data[5,45] = ""
$data[5,46] = ""
$data[5,47] = ""
$data[5,48] = 'Alexander Singer'
$data[5,49] = 'Alexander Singer'
$data[5,50] = ""

# Row 9
$data[6,0] = 101187646
$data[6,1] = 89410
$data[6,2] = 'Ovaliderad'
$data[6,3] = 'NT'
$data[6,4] = 5432
$data[6,5] = 'Granticka'
$data[6,6] = 'Porodaedalea chrysoloma'
$data[6,7] = '(Fr.) Fiasson & Niemelä'
$data[6,8] = ""
$data[6,9] = ""
$data[6,10] = ""
$data[6,11] = ""
$data[6,12] = ""
$data[6,13] = ""
$data[6,14] = ""
$data[6,15] = 'Knätten NB, Hjd'
$data[6,16] = 443331.2068214896
$data[6,17] = 6909321.897501401
$data[6,18] = 25
$data[6,19] = 'Jämtland'
$data[6,20] = 'Härjedalen'
$data[6,21] = 'Härjedalen'
$data[6,22] = 'Vemdalen'
$data[6,23] = ""
$data[6,24] = "'2022-05-27"
$data[6,25] = "'00:00"
$data[6,26] = "'2022-05-27"
$data[6,27] = "'00:00"
$data[6,28] = ""
$data[6,29] = $false
$data[6,30] = $false
$data[6,31] = ""
$data[6,32] = $false
$data[6,33] = ""
$data[6,34] = ""
$data[6,35] = ""
$data[6,36] = ""
$data[6,37] = ""
$data[6,38] = ""
$data[6,39] = ""
$data[6,40] = ""
$data[6,41] = ""
$data[6,42] = ""
$data[6,43] = ""
$data[6,44] = ""
$data[6,45] = ""
$data[6,46] = ""
$data[6,47] = ""
$data[6,48] = 'Alexander Singer'
$data[6,49] = 'Alexander Singer'
$data[6,50] = ""

# Row 10
$data[7,0] = 101188539
$data[7,1] = 96354
$data[7,2] = 'Ovaliderad'
$data[7,3] = 'LC'
$data[7,4] = 221952
$data[7,5] = 'Spindelblomster'
$data[7,6] = 'Neottia cordata'
$data[7,7] = '(L.) Rich.'
$data[7,8] = ""
$data[7,9] = ""
$data[7,10] = ""
$data[7,11] = ""
$data[7,12] = ""
$data[7,13] = ""
$data[7,14] = ""
$data[7,15] = 'Knätten NB, Hjd'
$data[7,16] = 443349.9363240626
$data[7,17] = 6909410.417016066
$data[7,18] = 25
$data[7,19] = 'Jämtland'
$data[7,20] = 'Härjedalen'
$data[7,21] = 'Härjedalen'
$data[7,22] = 'Vemdalen'
$data[7,23] = ""
$data[7,24] = "'2022-05-27"
$data[7,25] = "'00:00"
$data[7,26] = "'2022-05-27"
$data[7,27] = "'00:00"
$data[7,28] = ""
$data[7,29] = $false
$data[7,30] = $false
$data[7,31] = ""
$data[7,32] = $false
$data[7,33] = ""
$data[7,34] = ""
$data[7,35] = ""
$data[7,36] = ""
$data[7,37] = ""
$data[7,38] = ""
$data[7,39] = ""
$data[7,40] = ""
$data[7,41] = ""
$data[7,42] = ""
$data[7,43] = ""
$data[7,44] = ""
$data[7,45] = ""
$data[7,46] = ""
$data[7,47] = ""
$data[7,48] = 'Alexander Singer'
$data[7,49] = 'Alexander Singer'
$data[7,50] = ""

# Row 11
$data[8,0] = 101188150
$data[8,1] = 89338
$data[8,2] = 'Ovaliderad'
$data[8,3] = 'NT'
$data[8,4] = 112
$data[8,5] = 'Stjärntagging'
$data[8,6] = 'Asterodon ferruginosus'
$data[8,7] = 'Pat.'
$data[8,8] = ""
$data[8,9] = ""
$data[8,10] = ""
$data[8,11] = ""
$data[8,12] = ""
$data[8,13] = ""
$data[8,14] = ""
$data[8,15] = 'Knätten NB, Hjd'
$data[8,16] = 443320.66748062
$data[8,17] = 6909359.749755096
$data[8,18] = 25
$data[8,19] = 'Jämtland'
$data[8,20] = 'Härjedalen'
$data[8,21] = 'Härjedalen'
$data[8,22] = 'Vemdalen'
$data[8,23] = ""
$data[8,24] = "'2022-05-27"
$data[8,25] = "'00:00"
$data[8,26] = "'2022-05-27"
$data[8,27] = "'00:00"
$data[8,28] = ""
$data[8,29] = $false
$data[8,30] = $false
$data[8,31] = ""
$data[8,32] = $false
$data[8,33] = ""
$data[8,34] = ""
$data[8,35] = ""
$data[8,36] = ""
$data[8,37] = ""
$data[8,38] = ""
$data[8,39] = ""
$data[8,40] = ""
$data[8,41] = ""
$data[8,42] = ""
$data[8,43] = ""
$data[8,44] = ""
$data[8,45] = ""
$data[8,46] = ""
$data[8,47] = ""
$data[8,48] = 'Alexander Singer'
$data[8,49] = 'Alexander Singer'
$data[8,50] = ""

# Row 12
$data[9,0] = 101208904
$data[9,1] = 89392
$data[9,2] = 'Ovaliderad'
$data[9,3] = 'NT'
$data[9,4] = 1202
$data[9,5] = 'Ullticka'
$data[9,6] = 'Phellinidium ferrugineofuscum'
$data[9,7] = '(P.Karst.) Fiasson & Niemelä'
$data[9,8] = ""
$data[9,9] = ""
$data[9,10] = ""
$data[9,11] = ""
$data[9,12] = ""
$data[9,13] = ""
$data[9,14] = ""
$data[9,15] = 'Knätten NB, Hjd'
$data[9,16] = 443320.66748062
$data[9,17] = 6909359.749755096
$data[9,18] = 25
$data[9,19] = 'Jämtland'
$data[9,20] = 'Härjedalen'
$data[9,21] = 'Härjedalen'
$data[9,22] = 'Vemdalen'
$data[9,23] = ""
$data[9,24] = "'2022-05-27"
$data[9,25] = "'00:00"
$data[9,26] = "'2022-05-27"
$data[9,27] = "'00:00"
$data[9,28] = ""
$data[9,29] = $false
$data[9,30] = $false
$data[9,31] = ""
$data[9,32] = $false
$data[9,33] = ""
$data[9,34] = ""
$data[9,35] = ""
$data[9,36] = ""
$data[9,37] = ""
$data[9,38] = ""
$data[9,39] = ""
$data[9,40] = ""
$data[9,41] = ""
$data[9,42] = ""
$data[9,43] = ""
$data[9,44] = ""
$data[9,45] = ""
$data[9,46] = ""
$data[9,47] = ""
$data[9,48] = 'Alexander Singer'
$data[9,49] = 'Alexander Singer'
$data[9,50] = ""

# Row 13
$data[10,0] = 101208905
$data[10,1] = 89734
$data[10,2] = 'Ovaliderad'
$data[10,3] = 'VU'
$data[10,4] = 2063
$data[10,5] = 'Grantickeporing'
$data[10,6] = 'Skeletocutis chrysella'
$data[10,7] = 'Niemelä'
$data[10,8] = ""
$data[10,9] = ""
$data[10,10] = ""
$data[10,11] = ""
$data[10,12] = ""
$data[10,13] = ""
$data[10,14] = ""
$data[10,15] = 'Knätten NB, Hjd'
$data[10,16] = 443320.66748062
$data[10,17] = 6909359.749755096
$data[10,18] = 25
$data[10,19] = 'Jämtland'
$data[10,20] = 'Härjedalen'
$data[10,21] = 'Härjedalen'
$data[10,22] = 'Vemdalen'
$data[10,23] = ""
$data[10,24] = "'2022-05-27"
$data[10,25] = "'00:00"
$data[10,26] = "'2022-05-27"
$data[10,27] = "'00:00"
$data[10,28] = ""
$data[10,29] = $false
$data[10,30] = $false
$data[10,31] = ""
$data[10,32] = $false
$data[10,33] = ""
$data[10,34] = ""
$data[10,35] = ""
$data[10,36] = ""
$data[10,37] = ""
$data[10,38] = ""
$data[10,39] = ""
$data[10,40] = ""
$data[10,41] = ""
$data[10,42] = ""
$data[10,43] = ""
$data[10,44] = ""
$data[10,45] = ""
$data[10,46] = ""
$data[10,47] = ""
$data[10,48] = 'Alexander Singer'
$data[10,49] = 'Alexander Singer'
$data[10,50] = ""

$ws.Range("A3:AY13").Value = $data
Write-Output "Applied row permutation to A3:AY13"